$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each schedule date (column A) moves forward by exactly 1096 days
# (same month/day, year 2023 -> 2026). Column B holds the literal
# French weekday name for that date, so it is rewritten to match the
# new date's weekday.
#
# row -> (new date serial, new weekday name)
$updates = @(
    @{ Row = 2;  NewDate = 46056; Weekday = "mardi"    },
    @{ Row = 5;  NewDate = 46098; Weekday = "mardi"    },
    @{ Row = 8;  NewDate = 46104; Weekday = "lundi"    },
    @{ Row = 10; NewDate = 46105; Weekday = "mardi"    },
    @{ Row = 12; NewDate = 46106; Weekday = "mercredi" },
    @{ Row = 14; NewDate = 46108; Weekday = "vendredi" },
    @{ Row = 17; NewDate = 46112; Weekday = "mardi"    },
    @{ Row = 19; NewDate = 46114; Weekday = "jeudi"    },
    @{ Row = 22; NewDate = 46134; Weekday = "mercredi" },
    @{ Row = 26; NewDate = 46135; Weekday = "jeudi"    },
    @{ Row = 29; NewDate = 46136; Weekday = "vendredi" },
    @{ Row = 33; NewDate = 46147; Weekday = "mardi"    },
    @{ Row = 35; NewDate = 46149; Weekday = "jeudi"    },
    @{ Row = 39; NewDate = 46153; Weekday = "lundi"    },
    @{ Row = 41; NewDate = 46154; Weekday = "mardi"    },
    @{ Row = 44; NewDate = 46155; Weekday = "mercredi" }
)

foreach ($u in $updates) {
    $ws.Range("A" + $u.Row).Value2 = $u.NewDate
    $ws.Range("B" + $u.Row).Value = $u.Weekday
}
